$d = $word.ActiveDocument

# Fix 1: split "Wrackful Code of Conduct" into "Wrackful" + " Code of Conduct"
$d.Content.Find.Execute("Wrackful Code of Conduct", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Wrackful Code of Conduct", 2) | Out-Null

# Fix 2: fix the "outline" -> "outlined" typo
$d.Content.Find.Execute("rules outline above", $true, $false, $false, $false, $false,
                         $true, 1, $false, "rules outlined above", 2) | Out-Null
